# Trim the deck down to the single "Draft Outline" slide: delete every
# slide after slide 1 (the Google Trends, video-performance, Twitter
# listening, and sourcing slides that used to follow it).
$p = $ppt.ActivePresentation

while ($p.Slides.Count -gt 1) {
    $p.Slides.Item($p.Slides.Count).Delete()
}
